$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7:D12").Style = "Normal"

$ws.Range("A20").Value = "sc16"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = "A sublimit and a restriction on one of two policies"
$ws.Range("F20").Value = "complete"
$ws.Range("G20").Value = "yes"
$ws.Range("H20").Value = "done"
